# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Betarraga" ahead of the existing row 482, shifting the remaining rows
# (old 482-523) down by one to 483-524.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 482..523 down to 483..524, leaving row 482 free (and blank) for
# the new record. Excel.Rows.Insert shifts cells down and carries the
# formatting of the row above along with it, which is what keeps column D's
# date style intact.
$ws.Rows.Item(482).Insert()

# Populate the newly freed row 482 with the new weekly record.
$ws.Cells.Item(482, 1).Value = 4
$ws.Cells.Item(482, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(482, 3).Value = 'Los Lagos'
$ws.Cells.Item(482, 4).Value = 45166
$ws.Cells.Item(482, 5).Value = 10
$ws.Cells.Item(482, 6).Value = 100114014
$ws.Cells.Item(482, 7).Value = 'Betarraga'
$ws.Cells.Item(482, 8).Value = 'Sin especificar'
$ws.Cells.Item(482, 9).Value = 'Primera'
$ws.Cells.Item(482, 10).Value = 500
$ws.Cells.Item(482, 11).Value = 1000
$ws.Cells.Item(482, 12).Value = 1000
$ws.Cells.Item(482, 13).Value = 1000
$ws.Cells.Item(482, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(482, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(482, 16).Value = 200
$ws.Cells.Item(482, 17).Value = 5
$ws.Cells.Item(482, 18).Value = 'Hortaliza'
